$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 131 so the existing rows 131:151
# (Caqui records) shift down to 133:153, making room for a new pair
# of records (Mankaki Primera / Segunda) dated 44706.
$ws.Rows("131:132").Insert()

# New row 131 — Mankaki / Primera, fecha 44706
$ws.Cells.Item(131, 1).Value  = 6
$ws.Cells.Item(131, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(131, 3).Value  = "Metropolitana"
$ws.Cells.Item(131, 4).Value  = 44706
$ws.Cells.Item(131, 5).Value  = 13
$ws.Cells.Item(131, 6).Value  = "Fruta"
$ws.Cells.Item(131, 7).Value  = 100107
$ws.Cells.Item(131, 8).Value  = "Otros"
$ws.Cells.Item(131, 9).Value  = 100107001
$ws.Cells.Item(131, 10).Value = "Caqui"
$ws.Cells.Item(131, 11).Value = "Mankaki"
$ws.Cells.Item(131, 12).Value = "Primera"
$ws.Cells.Item(131, 13).Value = 14
$ws.Cells.Item(131, 14).Value = 280000
$ws.Cells.Item(131, 15).Value = 300000
$ws.Cells.Item(131, 16).Value = 290000
$ws.Cells.Item(131, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(131, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(131, 19).Value = 644
$ws.Cells.Item(131, 20).Value = 450

# New row 132 — Mankaki / Segunda, fecha 44706
$ws.Cells.Item(132, 1).Value  = 6
$ws.Cells.Item(132, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(132, 3).Value  = "Metropolitana"
$ws.Cells.Item(132, 4).Value  = 44706
$ws.Cells.Item(132, 5).Value  = 13
$ws.Cells.Item(132, 6).Value  = "Fruta"
$ws.Cells.Item(132, 7).Value  = 100107
$ws.Cells.Item(132, 8).Value  = "Otros"
$ws.Cells.Item(132, 9).Value  = 100107001
$ws.Cells.Item(132, 10).Value = "Caqui"
$ws.Cells.Item(132, 11).Value = "Mankaki"
$ws.Cells.Item(132, 12).Value = "Segunda"
$ws.Cells.Item(132, 13).Value = 12
$ws.Cells.Item(132, 14).Value = 240000
$ws.Cells.Item(132, 15).Value = 240000
$ws.Cells.Item(132, 16).Value = 240000
$ws.Cells.Item(132, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(132, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(132, 19).Value = 533
$ws.Cells.Item(132, 20).Value = 450

# Apply the date display format (matching the column's existing date cells)
# to the D cells of the two freshly inserted rows.
$ws.Range("D131:D132").NumberFormat = "YYYY-MM-DD HH:MM:SS"
